$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1958041958041958
$ws.Range("C2").Value = 0.5571095571095571
$ws.Range("J2").Value = 0.02447552447552448
$ws.Range("P2").Value = 0.1351981351981352
$ws.Range("S2").Value = 0.08741258741258741
$ws.Range("B3").Value = 0.009940357852882704
$ws.Range("C3").Value = 0.03180914512922465
$ws.Range("J3").Value = 0.02783300198807157
$ws.Range("P3").Value = 0.7435387673956262
$ws.Range("S3").Value = 0.1868787276341948
$ws.Range("J4").Value = 0.03759398496240601
$ws.Range("P4").Value = 0.7443609022556391
$ws.Range("S4").Value = 0.2180451127819549
$ws.Range("D6").Value = 0.015625
$ws.Range("F6").Value = 0.08159722222222222
$ws.Range("J6").Value = 0.2413194444444444
$ws.Range("O6").Value = 0.01736111111111111
$ws.Range("Q6").Value = 0.1614583333333333
$ws.Range("R6").Value = 0.06076388888888889
$ws.Range("B7").Value = 0.1043165467625899
$ws.Range("D7").Value = 0.02697841726618705
$ws.Range("E7").Value = 0.001798561151079137
$ws.Range("F7").Value = 0.05575539568345324
$ws.Range("J7").Value = 0.1384892086330935
$ws.Range("O7").Value = 0.02338129496402878
$ws.Range("Q7").Value = 0.1924460431654676
$ws.Range("R7").Value = 0.06654676258992806
$ws.Range("S7").Value = 0.3902877697841727
$ws.Range("B8").Value = 0.09856781802864364
$ws.Range("D8").Value = 0.0134793597304128
$ws.Range("F8").Value = 0.05560235888795282
$ws.Range("J8").Value = 0.1112047177759056
$ws.Range("O8").Value = 0.02274641954507161
$ws.Range("Q8").Value = 0.1760741364785173
$ws.Range("R8").Value = 0.1061499578770008
$ws.Range("S8").Value = 0.4161752316764954
$ws.Range("B9").Value = 0.09913043478260869
$ws.Range("D9").Value = 0.02608695652173913
$ws.Range("E9").Value = 0.001739130434782609
$ws.Range("F9").Value = 0.05043478260869565
$ws.Range("J9").Value = 0.1478260869565217
$ws.Range("O9").Value = 0.01565217391304348
$ws.Range("Q9").Value = 0.1686956521739131
$ws.Range("R9").Value = 0.07130434782608695
$ws.Range("S9").Value = 0.4191304347826087
$ws.Range("B10").Value = 0.1103678929765886
$ws.Range("D10").Value = 0.02173913043478261
$ws.Range("E10").Value = 0.0005574136008918618
$ws.Range("F10").Value = 0.05964325529542921
$ws.Range("J10").Value = 0.141025641025641
$ws.Range("O10").Value = 0.01700111482720178
$ws.Range("Q10").Value = 0.2159977703455964
$ws.Range("R10").Value = 0.0794314381270903
$ws.Range("S10").Value = 0.3542363433667782
$ws.Range("G11").Value = 0.1457831325301205
$ws.Range("J11").Value = 0.0819277108433735
$ws.Range("K11").Value = 0.1939759036144578
$ws.Range("L11").Value = 0.5698795180722892
$ws.Range("S11").Value = 0.008433734939759036
$ws.Range("F12").Value = 0.002028397565922921
$ws.Range("G12").Value = 0.742393509127789
$ws.Range("J12").Value = 0.1825557809330629
$ws.Range("K12").Value = 0.006085192697768763
$ws.Range("L12").Value = 0.03651115618661258
$ws.Range("S12").Value = 0.03042596348884381
$ws.Range("F13").Value = 0.008
$ws.Range("G13").Value = 0.672
$ws.Range("J13").Value = 0.264
$ws.Range("S13").Value = 0.056
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.02360876897133221
$ws.Range("H15").Value = 0.1281618887015177
$ws.Range("I15").Value = 0.06913996627318718
$ws.Range("J15").Value = 0.3777403035413153
$ws.Range("K15").Value = 0.07419898819561552
$ws.Range("M15").Value = 0.01854974704890388
$ws.Range("O15").Value = 0.07419898819561552
$ws.Range("S15").Value = 0.2344013490725126
$ws.Range("F16").Value = 0.01601423487544484
$ws.Range("H16").Value = 0.1903914590747331
$ws.Range("I16").Value = 0.07651245551601424
$ws.Range("J16").Value = 0.4270462633451957
$ws.Range("K16").Value = 0.103202846975089
$ws.Range("M16").Value = 0.0195729537366548
$ws.Range("N16").Value = 0.001779359430604982
$ws.Range("O16").Value = 0.05160142348754448
$ws.Range("S16").Value = 0.1138790035587189
$ws.Range("F17").Value = 0.01888276947285602
$ws.Range("H17").Value = 0.1730920535011802
$ws.Range("I17").Value = 0.09992132179386309
$ws.Range("J17").Value = 0.3949645948072384
$ws.Range("K17").Value = 0.09992132179386309
$ws.Range("M17").Value = 0.01730920535011802
$ws.Range("N17").Value = 0.0007867820613690008
$ws.Range("O17").Value = 0.07317073170731707
$ws.Range("S17").Value = 0.1219512195121951
$ws.Range("F18").Value = 0.02651515151515152
$ws.Range("H18").Value = 0.2045454545454546
$ws.Range("I18").Value = 0.07007575757575757
$ws.Range("J18").Value = 0.4337121212121212
$ws.Range("K18").Value = 0.0928030303030303
$ws.Range("M18").Value = 0.02272727272727273
$ws.Range("N18").Value = 0.001893939393939394
$ws.Range("O18").Value = 0.04166666666666666
$ws.Range("S18").Value = 0.1060606060606061
$ws.Range("F19").Value = 0.01748251748251748
$ws.Range("H19").Value = 0.1963869463869464
$ws.Range("I19").Value = 0.0944055944055944
$ws.Range("J19").Value = 0.3744172494172494
$ws.Range("K19").Value = 0.1139277389277389
$ws.Range("M19").Value = 0.02068764568764569
$ws.Range("N19").Value = 0.0002913752913752914
$ws.Range("O19").Value = 0.06177156177156177
$ws.Range("S19").Value = 0.1206293706293706
